# Fix the "Last Review Date" column: refresh the dates (stored as serial
# numbers, same as Excel does internally) and switch the display format
# from dd/mm/yy to m/d/yyyy. Also pre-format a block of rows below the
# existing data (F7:F21) so newly uploaded rows already carry the date
# format (NGRX store now appends rows there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date values (serial numbers) for the existing data rows 2-6, column F
$dates = @(45667, 45668, 45669, 45670, 45671)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $dates[$i]
}

# Apply the updated date number format to F2:F21 (covers the existing rows
# plus the freshly pre-formatted, still-empty rows 7:21)
$ws.Range("F2:F21").NumberFormat = "m/d/yyyy;@"

# Move selection to F2, matching the post-edit selection in the workbook
$ws.Range("F2").Select()
